# Regenerate Report for Archive
#
# The localization-status workbook tracks, per source file, its current
# handoff/handback status. Two rows (the 95be4f59...md and
# b96d7e22...md entries) need to trade places - b96d7e22 moves up to
# row 4 (now "In Translation") and 95be4f59 drops to row 5 (still
# "Ready for handoff") - on every sheet: the Overview summary plus the
# per-locale zh-cn and de-de detail sheets. The hyperlinks already
# attached to cells A4/A5 (and C4/C5 on the locale sheets) point at the
# correct per-file URLs and must stay put; only the displayed text
# moves.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (columns A:C) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A4").Value = "b96d7e22-f461-48af-95b0-d05193254f94.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("A5").Value = "95be4f59-b5c5-47e2-bbde-e978b4234152.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

# --- zh-cn detail sheet (columns A:D) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A4").Value = "b96d7e22-f461-48af-95b0-d05193254f94.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "b96d7e22-f461-48af-95b0-d05193254f94.098a529e1403db042deeceefeb618ac8b7419cb2.zh-cn.xlf"
$ws.Range("D4").Value = "2016-02-17 04:03:30"
$ws.Range("A5").Value = "95be4f59-b5c5-47e2-bbde-e978b4234152.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "95be4f59-b5c5-47e2-bbde-e978b4234152.8eaa6fa3d63e10aeab935d719f0fd29648066b5b.zh-cn.xlf"
$ws.Range("D5").Value = "2016-02-17 04:04:12"

# --- de-de detail sheet (columns A:D) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A4").Value = "b96d7e22-f461-48af-95b0-d05193254f94.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "b96d7e22-f461-48af-95b0-d05193254f94.098a529e1403db042deeceefeb618ac8b7419cb2.de-de.xlf"
$ws.Range("D4").Value = "2016-02-17 04:03:40"
$ws.Range("A5").Value = "95be4f59-b5c5-47e2-bbde-e978b4234152.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "95be4f59-b5c5-47e2-bbde-e978b4234152.8eaa6fa3d63e10aeab935d719f0fd29648066b5b.de-de.xlf"
$ws.Range("D5").Value = "2016-02-17 04:04:21"
